$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 18 (pushes e501.. rows down by one)
$ws.Rows("18").Insert()

# New e017 entry
$ws.Range("A18").Value = "e017"

$e017Body = @'
<Bold>e017 Preparations Final</Bold> 
<InlineUIContainer><Button Content='r4.46' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>    
<LineBreak/><LineBreak/>
US Control markers are placed on sectors 1, 2, and 3. The Weather is displayed on top left of Battle Board.
<LineBreak/><LineBreak/>
                                            <InlineUIContainer><Image Name='Continue017' Height='100' Width='100'></Image></InlineUIContainer>
'@
$ws.Range("B18").Value = $e017Body

# Row height adjustments
$ws.Rows("5").RowHeight = 105
$ws.Rows("6").RowHeight = 120
$ws.Rows("18").RowHeight = 90

# Update the current selection to match the final view
$ws.Range("B21").Select()
